# Horarios actualizados Linea 141 - 49
# Insert the newly-scraped row (02:29:13 / 04:01 / 81_EL PELIGRO / 92) ahead
# of the existing rows on both the "LP1912" and "6203-6173" sheets, bump the
# "Ultima actualizacion" timestamp + "Total filas" count, and leave the
# "LP1912-215" sheet untouched (it has no changes in this update).

$wb = $excel.ActiveWorkbook

$sheetNames = @("LP1912", "6203-6173")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A2").Value = "Última actualización: 02:29:13"
    $ws.Range("A3").Value = "Total filas: 6"

    # Push the existing data rows (old row 9 -> 10, old row 10 -> 11) down to
    # make room for the new row at position 9.
    $ws.Rows.Item(9).Insert()

    $ws.Range("A9").Value = "02:29:13"
    $ws.Range("B9").Value = "04:01"
    $ws.Range("C9").Value = "81_EL PELIGRO"
    $ws.Range("D9").Value = 92
}
